# Atualização de bases das ligas, do dia: 17-06-2024 às 21:10
# Two pairs of adjacent match rows had been recorded in the wrong order;
# this swaps the full record (every column except the running index in A,
# which stays tied to its row) between each pair so the rows land on the
# correct match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-MatchRows {
    param($ws, $row1, $row2)

    $range1 = $ws.Range("B$row1" + ":AD$row1")
    $range2 = $ws.Range("B$row2" + ":AD$row2")

    $values1 = $range1.Value2
    $values2 = $range2.Value2

    $range1.Value2 = $values2
    $range2.Value2 = $values1
}

Swap-MatchRows $ws 107 108
Swap-MatchRows $ws 128 129
Swap-MatchRows $ws 211 212
Swap-MatchRows $ws 214 215
